$d = $word.ActiveDocument

# Remove the trailing "Requisitos" section: its Heading2 title paragraph
# and the following ListBullet paragraph listing the prerequisite
# "LOB1046 - Engenharia do Meio Ambiente (Requisito fraco)". These are
# the last two paragraphs of the document, immediately after the
# Bibliografia content paragraph.

$count = $d.Paragraphs.Count
$last = $d.Paragraphs($count)
$secondLast = $d.Paragraphs($count - 1)

$deleteRange = $d.Range($secondLast.Range.Start, $last.Range.End)
$deleteRange.Delete()
